$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.440.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.574.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +1.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.491'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.17'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '23.74'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.22%  '
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0880'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.798.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.571.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.18%  '
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.417.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("E21").Value = '  -1.60%  '
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("E23").Value = '  -4.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("E25").Value = '  +3.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("E29").Value = '  -2.32%  '
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("E31").Value = '  -2.03%  '
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.391.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("E36").Value = '  -1.97%  '
$ws.Range("E37").Value = '  -3.39%  '
$ws.Range("E38").Value = '  +3.13%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0165'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.55'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.44%  '
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.792'
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.978'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.31'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.711.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("E50").Value = '  +2.33%  '
$ws.Range("E51").Value = '  -1.27%  '
